$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 (bold, bordered) into F1, then set its value.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
  "2021-10-05 10:52:00.978488",
  "2021-10-05 10:52:00.978501",
  "2021-10-05 10:52:00.978504",
  "2021-10-05 10:52:00.978508",
  "2021-10-05 10:52:00.978511",
  "2021-10-05 10:52:00.978514",
  "2021-10-05 10:52:00.978518",
  "2021-10-05 10:52:00.978521",
  "2021-10-05 10:52:00.978524",
  "2021-10-05 10:52:00.978528",
  "2021-10-05 10:52:00.978531",
  "2021-10-05 10:52:00.978534",
  "2021-10-05 10:52:00.978537",
  "2021-10-05 10:52:00.978540",
  "2021-10-05 10:52:00.978543",
  "2021-10-05 10:52:00.978546",
  "2021-10-05 10:52:00.978550",
  "2021-10-05 10:52:00.978553",
  "2021-10-05 10:52:00.978556",
  "2021-10-05 10:52:00.978559",
  "2021-10-05 10:52:00.978562",
  "2021-10-05 10:52:00.978565",
  "2021-10-05 10:52:00.978568",
  "2021-10-05 10:52:00.978571",
  "2021-10-05 10:52:00.978575",
  "2021-10-05 10:52:00.978578",
  "2021-10-05 10:52:00.978581",
  "2021-10-05 10:52:00.978584",
  "2021-10-05 10:52:00.978587",
  "2021-10-05 10:52:00.978590",
  "2021-10-05 10:52:00.978593",
  "2021-10-05 10:52:00.978596",
  "2021-10-05 10:52:00.978600",
  "2021-10-05 10:52:00.978603",
  "2021-10-05 10:52:00.978606",
  "2021-10-05 10:52:00.978609",
  "2021-10-05 10:52:00.978612"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
  $row = $i + 2
  $cell = $ws.Cells.Item($row, 6)
  $cell.Value = $timestamps[$i]
}
